# Correção das notas do fórum para matc65 em 2021.2
# For every student row whose "nota_view" (column J) equals 4, zero out
# all the daily view counters (columns B:H) plus the computed totals
# (total_views in I, nota_view in J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $nota = $ws.Cells.Item($r, 10).Value2  # column J = nota_view
    if ($nota -eq 4) {
        $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 10)).Value2 = 0
    }
}
